$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 452.7742
$ws.Cells.Item(19, 9).Value = 350.5
$ws.Cells.Item(19, 10).Value = 537
$ws.Cells.Item(19, 11).Value = 350.5
$ws.Cells.Item(19, 12).Value = 537
$ws.Cells.Item(19, 13).Value = -175.5
$ws.Cells.Item(19, 14).Value = -887

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 2417.2
$ws.Cells.Item(51, 9).Value = 2780
$ws.Cells.Item(51, 11).Value = 2780
$ws.Cells.Item(51, 13).Value = -2296

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 3775
$ws.Cells.Item(70, 9).Value = 900
$ws.Cells.Item(70, 10).Value = 4733.3335
$ws.Cells.Item(70, 11).Value = 2700
$ws.Cells.Item(70, 12).Value = 14200.0005
$ws.Cells.Item(70, 13).Value = -2430
$ws.Cells.Item(70, 14).Value = -14740.0005

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 3775
$ws.Cells.Item(73, 9).Value = 900
$ws.Cells.Item(73, 10).Value = 4733.3335
$ws.Cells.Item(73, 11).Value = 2700
$ws.Cells.Item(73, 12).Value = 14200.0005
$ws.Cells.Item(73, 13).Value = -1764
$ws.Cells.Item(73, 14).Value = -16072.0005

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 42311384
$ws.Cells.Item(76, 9).Value = 44003520
$ws.Cells.Item(76, 11).Value = 44003520
$ws.Cells.Item(76, 13).Value = -44003205

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(79, 8).Value = 42311384
$ws.Cells.Item(79, 9).Value = 44003520
$ws.Cells.Item(79, 11).Value = 44003520
$ws.Cells.Item(79, 13).Value = -44002428

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 409984.53
$ws.Cells.Item(80, 9).Value = 758.5
$ws.Cells.Item(80, 10).Value = 1555817.4
$ws.Cells.Item(80, 11).Value = 2275.5
$ws.Cells.Item(80, 12).Value = 4667452.199999999
$ws.Cells.Item(80, 13).Value = -1277.5
$ws.Cells.Item(80, 14).Value = -4669448.199999999

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(83, 8).Value = 409984.53
$ws.Cells.Item(83, 9).Value = 758.5
$ws.Cells.Item(83, 10).Value = 1555817.4
$ws.Cells.Item(83, 11).Value = 6826.5
$ws.Cells.Item(83, 12).Value = 14002356.6
$ws.Cells.Item(83, 13).Value = -1834.5
$ws.Cells.Item(83, 14).Value = -14012340.6

# ALC row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(96, 8).Value = 4040.75
$ws.Cells.Item(96, 9).Value = 4699.9375
$ws.Cells.Item(96, 10).Value = 2722.375
$ws.Cells.Item(96, 11).Value = 14099.8125
$ws.Cells.Item(96, 12).Value = 8167.125
$ws.Cells.Item(96, 13).Value = -12726.8125
$ws.Cells.Item(96, 14).Value = -10913.125

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1050
$ws.Cells.Item(98, 9).Value = 988.63635
$ws.Cells.Item(98, 10).Value = 1500
$ws.Cells.Item(98, 11).Value = 988.63635
$ws.Cells.Item(98, 12).Value = 1500
$ws.Cells.Item(98, 13).Value = 509.36365
$ws.Cells.Item(98, 14).Value = -4496

# ALC row 121
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(121, 8).Value = 867.2174
$ws.Cells.Item(121, 10).Value = 873.619
$ws.Cells.Item(121, 12).Value = 2620.857
$ws.Cells.Item(121, 14).Value = -6114.857

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 1050
$ws.Cells.Item(122, 9).Value = 988.63635
$ws.Cells.Item(122, 10).Value = 1500
$ws.Cells.Item(122, 11).Value = 2965.90905
$ws.Cells.Item(122, 12).Value = 4500
$ws.Cells.Item(122, 13).Value = -515.9090500000002
$ws.Cells.Item(122, 14).Value = -9400

# ALC row 124
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(124, 8).Value = 39000
$ws.Cells.Item(124, 10).Value = 39000
$ws.Cells.Item(124, 12).Value = 39000
$ws.Cells.Item(124, 14).Value = -48820

# ALC row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(127, 8).Value = 2291.0227
$ws.Cells.Item(127, 9).Value = 866.1667
$ws.Cells.Item(127, 10).Value = 2516
$ws.Cells.Item(127, 11).Value = 2598.5001
$ws.Cells.Item(127, 12).Value = 7548
$ws.Cells.Item(127, 13).Value = 2361.4999
$ws.Cells.Item(127, 14).Value = -17468

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 11113329
$ws.Cells.Item(137, 9).Value = 1871.4348
$ws.Cells.Item(137, 10).Value = 22729854
$ws.Cells.Item(137, 11).Value = 5614.3044
$ws.Cells.Item(137, 12).Value = 68189562
$ws.Cells.Item(137, 13).Value = -3064.3044
$ws.Cells.Item(137, 14).Value = -68194662

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 20173.902
$ws.Cells.Item(32, 9).Value = 18919.836
$ws.Cells.Item(32, 10).Value = 30345.777
$ws.Cells.Item(32, 11).Value = 18919.836
$ws.Cells.Item(32, 12).Value = 30345.777
$ws.Cells.Item(32, 13).Value = -18632.836
$ws.Cells.Item(32, 14).Value = -30919.777

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 2645.5625
$ws.Cells.Item(88, 9).Value = 2488.6
$ws.Cells.Item(88, 10).Value = 5000
$ws.Cells.Item(88, 11).Value = 2488.6
$ws.Cells.Item(88, 12).Value = 5000
$ws.Cells.Item(88, 13).Value = -2082.6
$ws.Cells.Item(88, 14).Value = -5812

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 2645.5625
$ws.Cells.Item(91, 9).Value = 2488.6
$ws.Cells.Item(91, 10).Value = 5000
$ws.Cells.Item(91, 11).Value = 2488.6
$ws.Cells.Item(91, 12).Value = 5000
$ws.Cells.Item(91, 13).Value = -1084.6
$ws.Cells.Item(91, 14).Value = -7808

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 8801.666999999999
$ws.Cells.Item(105, 9).Value = 4562
$ws.Cells.Item(105, 11).Value = 4562
$ws.Cells.Item(105, 13).Value = -2815

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 27660
$ws.Cells.Item(107, 9).Value = 1166.6666
$ws.Cells.Item(107, 10).Value = 67400
$ws.Cells.Item(107, 11).Value = 1166.6666
$ws.Cells.Item(107, 12).Value = 67400
$ws.Cells.Item(107, 13).Value = 753.3334
$ws.Cells.Item(107, 14).Value = -71240

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 159.52942
$ws.Cells.Item(7, 9).Value = 144.28572
$ws.Cells.Item(7, 10).Value = 230.66667
$ws.Cells.Item(7, 11).Value = 144.28572
$ws.Cells.Item(7, 12).Value = 230.66667
$ws.Cells.Item(7, 13).Value = -31.28572
$ws.Cells.Item(7, 14).Value = -456.66667

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1881.9149
$ws.Cells.Item(31, 9).Value = 1076.6923
$ws.Cells.Item(31, 10).Value = 2189.7942
$ws.Cells.Item(31, 11).Value = 1076.6923
$ws.Cells.Item(31, 12).Value = 2189.7942
$ws.Cells.Item(31, 13).Value = -781.6922999999999
$ws.Cells.Item(31, 14).Value = -2779.7942

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 1881.9149
$ws.Cells.Item(34, 9).Value = 1076.6923
$ws.Cells.Item(34, 10).Value = 2189.7942
$ws.Cells.Item(34, 11).Value = 1076.6923
$ws.Cells.Item(34, 12).Value = 2189.7942
$ws.Cells.Item(34, 13).Value = -874.6922999999999
$ws.Cells.Item(34, 14).Value = -2593.7942

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 886.65
$ws.Cells.Item(68, 9).Value = 505.93442
$ws.Cells.Item(68, 11).Value = 1517.80326
$ws.Cells.Item(68, 13).Value = -706.8032599999999

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(71, 8).Value = 886.65
$ws.Cells.Item(71, 9).Value = 505.93442
$ws.Cells.Item(71, 11).Value = 4553.40978
$ws.Cells.Item(71, 13).Value = -497.40978

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4386.364
$ws.Cells.Item(70, 9).Value = 4105.5557
$ws.Cells.Item(70, 10).Value = 5650
$ws.Cells.Item(70, 11).Value = 4105.5557
$ws.Cells.Item(70, 12).Value = 5650
$ws.Cells.Item(70, 13).Value = -3835.5557
$ws.Cells.Item(70, 14).Value = -6190

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 4386.364
$ws.Cells.Item(73, 9).Value = 4105.5557
$ws.Cells.Item(73, 10).Value = 5650
$ws.Cells.Item(73, 11).Value = 4105.5557
$ws.Cells.Item(73, 12).Value = 5650
$ws.Cells.Item(73, 13).Value = -3169.5557
$ws.Cells.Item(73, 14).Value = -7522

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3197.125
$ws.Cells.Item(132, 9).Value = 2305
$ws.Cells.Item(132, 10).Value = 5873.5
$ws.Cells.Item(132, 11).Value = 6915
$ws.Cells.Item(132, 12).Value = 17620.5
$ws.Cells.Item(132, 13).Value = -4385
$ws.Cells.Item(132, 14).Value = -22680.5

# GSM row 138
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(138, 8).Value = 39835
$ws.Cells.Item(138, 10).Value = 39835
$ws.Cells.Item(138, 12).Value = 39835
$ws.Cells.Item(138, 14).Value = -50115

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 45834.668
$ws.Cells.Item(40, 9).Value = 50001.6
$ws.Cells.Item(40, 10).Value = 25000
$ws.Cells.Item(40, 11).Value = 50001.6
$ws.Cells.Item(40, 12).Value = 25000
$ws.Cells.Item(40, 13).Value = -49865.6
$ws.Cells.Item(40, 14).Value = -25272

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 4369.1763
$ws.Cells.Item(122, 9).Value = 4938.731
$ws.Cells.Item(122, 10).Value = 2518.125
$ws.Cells.Item(122, 11).Value = 14816.193
$ws.Cells.Item(122, 12).Value = 7554.375
$ws.Cells.Item(122, 13).Value = -12366.193
$ws.Cells.Item(122, 14).Value = -12454.375

# LTW row 134
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(134, 8).Value = 67557.25
$ws.Cells.Item(134, 10).Value = 67557.25
$ws.Cells.Item(134, 12).Value = 67557.25
$ws.Cells.Item(134, 14).Value = -77697.25

# LTW row 139
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(139, 8).Value = 49703.75
$ws.Cells.Item(139, 10).Value = 49703.75
$ws.Cells.Item(139, 12).Value = 49703.75
$ws.Cells.Item(139, 14).Value = -59983.75

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 2903.4
$ws.Cells.Item(96, 10).Value = 2936.8
$ws.Cells.Item(96, 12).Value = 2936.8
$ws.Cells.Item(96, 14).Value = -5682.8

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 782
$ws.Cells.Item(113, 9).Value = 740.4
$ws.Cells.Item(113, 10).Value = 990
$ws.Cells.Item(113, 11).Value = 2221.2
$ws.Cells.Item(113, 12).Value = 2970
$ws.Cells.Item(113, 13).Value = -51.19999999999982
$ws.Cells.Item(113, 14).Value = -7310

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1552.9412
$ws.Cells.Item(122, 9).Value = 1150
$ws.Cells.Item(122, 11).Value = 3450
$ws.Cells.Item(122, 13).Value = -1000

# WVR row 125
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(125, 8).Value = 59000
$ws.Cells.Item(125, 10).Value = 59000
$ws.Cells.Item(125, 12).Value = 59000
$ws.Cells.Item(125, 14).Value = -68840

# WVR row 133
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(133, 8).Value = 68778.75
$ws.Cells.Item(133, 10).Value = 68778.75
$ws.Cells.Item(133, 12).Value = 68778.75
$ws.Cells.Item(133, 14).Value = -78898.75

# WVR row 138
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(138, 8).Value = 47680
$ws.Cells.Item(138, 10).Value = 47680
$ws.Cells.Item(138, 12).Value = 47680
$ws.Cells.Item(138, 14).Value = -57960
